# Nexial "#system" reference sheet: add the new `outputToCloud(resource)`
# command to the `base` module, and add a brand-new `text` module whose
# sole command is `spellCheck(var,profile,text)`.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1. Column E ("base" module): insert "outputToCloud(resource)" in its
#    alphabetically-sorted position (between "macro(file,sheet,name)" at
#    E21 and "prependText(var,prependWith)" at E22), pushing E22:E38 down
#    to E23:E39.
# ---------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $ws.Range("E" + ($r + 1)).Value2 = $ws.Range("E" + $r).Value2
}
$ws.Range("E22").Value2 = "outputToCloud(resource)"

# ---------------------------------------------------------------------
# 2. Column A ("target" list of module names): insert "text" in its
#    alphabetically-sorted position (between "step" at A24 and "web" at
#    A25), pushing A25:A30 down to A26:A31.
# ---------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $ws.Range("A" + ($r + 1)).Value2 = $ws.Range("A" + $r).Value2
}
$ws.Range("A25").Value2 = "text"

# ---------------------------------------------------------------------
# 3. Insert a brand-new column at Y (shifting web/webalert/webcookie/
#    ws/ws.async/xml from Y..AD to Z..AE) and populate the new Y column
#    with the new "text" module: header in row 1, single command in
#    row 2.
# ---------------------------------------------------------------------
$ws.Range("Y:Y").Insert()
$ws.Range("Y1").Value2 = "text"
$ws.Range("Y2").Value2 = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------
# 4. Fix up the defined names so they reference the correct (shifted)
#    ranges, and register the new "text" name.
# ---------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
